$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Journal de travail - new entries (rows 6-13), continuing after the
# existing "Découverte du projet..." entry on row 5.
# Dates are written as raw date serial numbers (matching the workbook's
# existing 1900 date system) so there is no ambiguity from locale-specific
# string parsing; the cells already carry the MM/DD/YY date style.

$ws.Range("A6").Value = 43001
$ws.Range("B6").Value = "BD : Modèle conceptuel"
$ws.Range("C6").Value = 3

$ws.Range("A7").Value = 43002
$ws.Range("B7").Value = "Liste des fonctionnalités"
$ws.Range("C7").Value = 3

$ws.Range("A8").Value = 43009
$ws.Range("B8").Value = "Organisation du serveur"
$ws.Range("C8").Value = 2.5

$ws.Range("A9").Value = 43011
$ws.Range("B9").Value = "Diagramme de classes"
$ws.Range("C9").Value = 3.5

$ws.Range("A10").Value = 43011
$ws.Range("B10").Value = "Paramétrage de l’ORM"
$ws.Range("C10").Value = 1.5

$ws.Range("A11").Value = 43017
$ws.Range("B11").Value = "Architecture client-serveur"
$ws.Range("C11").Value = 4

$ws.Range("A12").Value = 43018
$ws.Range("B12").Value = "Révision diagramme de classes"
$ws.Range("C12").Value = 1

$ws.Range("A13").Value = 43018
$ws.Range("B13").Value = "BD : Modèle logique"
$ws.Range("C13").Value = 2.5

# C32 already holds =SUM(C5:C31), so the total recalculates automatically.

$ws.Range("A14").Select()
